# Apply the diff:
#   - swap the data rows for (7,8) and (12,13)  [columns A:N]
#   - refresh the scrape timestamp in column O for every data row (2-37)
#
# Columns A,B,C,D,G,H,I,J,K,L,M,N hold text (even when they look numeric,
# e.g. "2.30"), while columns E and F hold real numbers. We capture both
# a Text (string) snapshot and a Value2 (native) snapshot of each row so
# the swap can restore the correct type per column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCols = @(1,2,3,4,7,8,9,10,11,12,13,14)   # A,B,C,D,G,H,I,J,K,L,M,N
$numCols  = @(5,6)                            # E,F

function Get-RowSnapshot($row) {
    $snap = @{}
    foreach ($c in $textCols) {
        $snap[$c] = $ws.Cells.Item($row, $c).Text
    }
    foreach ($c in $numCols) {
        $cell = $ws.Cells.Item($row, $c)
        if ($cell.Text -eq "") {
            $snap[$c] = $null
        } else {
            $snap[$c] = $cell.Value2
        }
    }
    return $snap
}

function Set-RowSnapshot($row, $snap) {
    foreach ($c in $textCols) {
        $cell = $ws.Cells.Item($row, $c)
        $text = $snap[$c]
        if ($text -eq "") {
            $cell.Value = ""
        } else {
            # Leading apostrophe forces literal-text storage (quotePrefix),
            # so numeric-looking strings like "2.30" keep their exact
            # formatting instead of being parsed into a float.
            $cell.Value = "'" + $text
        }
    }
    foreach ($c in $numCols) {
        $cell = $ws.Cells.Item($row, $c)
        if ($null -eq $snap[$c]) {
            $cell.Value = $null
        } else {
            $cell.Value = $snap[$c]
        }
    }
}

function Swap-Rows($r1, $r2) {
    $snap1 = Get-RowSnapshot $r1
    $snap2 = Get-RowSnapshot $r2
    Set-RowSnapshot $r1 $snap2
    Set-RowSnapshot $r2 $snap1
}

# --- Swap row 7 <-> row 8 ---
Swap-Rows 7 8

# --- Swap row 12 <-> row 13 ---
Swap-Rows 12 13

# --- Refresh timestamp column O for every data row (2-37) ---
$newTimestamp = "2022-07-24 20:58:22"
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 15).Value = "'" + $newTimestamp
}
